$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7167.5
$ws.Range("I62").Value = 4500.6665
$ws.Range("K62").Value = 4500.6665
$ws.Range("M62").Value = -3876.6665
$ws.Range("H65").Value = 7167.5
$ws.Range("I65").Value = 4500.6665
$ws.Range("K65").Value = 22503.3325
$ws.Range("M65").Value = -19383.3325
$ws.Range("H94").Value = 2490.4
$ws.Range("J94").Value = 500
$ws.Range("L94").Value = 500
$ws.Range("N94").Value = -1402
$ws.Range("H103").Value = 1339.3
$ws.Range("J103").Value = 1978.3334
$ws.Range("L103").Value = 5935.0002
$ws.Range("N103").Value = -7107.0002
$ws.Range("H127").Value = 905.6667
$ws.Range("J127").Value = 3000
$ws.Range("L127").Value = 9000
$ws.Range("N127").Value = -18920
$ws.Range("H137").Value = 3010.5386
$ws.Range("I137").Value = 2472.6365
$ws.Range("J137").Value = 3120.111
$ws.Range("K137").Value = 7417.9095
$ws.Range("L137").Value = 9360.332999999999
$ws.Range("M137").Value = -4867.9095
$ws.Range("N137").Value = -14460.333
$ws.Range("H138").Value = 3298.8816
$ws.Range("I138").Value = 2710.7273
$ws.Range("J138").Value = 3398.4153
$ws.Range("K138").Value = 8132.1819
$ws.Range("L138").Value = 10195.2459
$ws.Range("M138").Value = -2992.1819
$ws.Range("N138").Value = -20475.2459
$ws.Range("H141").Value = 2792.318
$ws.Range("I141").Value = 2707.2104
$ws.Range("J141").Value = 3331.3333
$ws.Range("K141").Value = 8121.6312
$ws.Range("L141").Value = 9993.999899999999
$ws.Range("M141").Value = -2941.6312
$ws.Range("N141").Value = -20353.9999

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6934.048
$ws.Range("I2").Value = 1156.7646
$ws.Range("J2").Value = 31487.5
$ws.Range("K2").Value = 1156.7646
$ws.Range("L2").Value = 31487.5
$ws.Range("M2").Value = -1043.7646
$ws.Range("N2").Value = -31713.5
$ws.Range("H32").Value = 2977.0317
$ws.Range("I32").Value = 2418.951
$ws.Range("J32").Value = 19998.5
$ws.Range("K32").Value = 2418.951
$ws.Range("L32").Value = 19998.5
$ws.Range("M32").Value = -2131.951
$ws.Range("N32").Value = -20572.5
$ws.Range("H61").Value = 4124.1304
$ws.Range("I61").Value = 3354.182
$ws.Range("J61").Value = 4829.9165
$ws.Range("K61").Value = 3354.182
$ws.Range("L61").Value = 4829.9165
$ws.Range("M61").Value = -3142.182
$ws.Range("N61").Value = -5253.9165
$ws.Range("H74").Value = 13891408
$ws.Range("I74").Value = 15153249
$ws.Range("K74").Value = 15153249
$ws.Range("M74").Value = -15152375
$ws.Range("H77").Value = 13891408
$ws.Range("I77").Value = 15153249
$ws.Range("K77").Value = 75766245
$ws.Range("M77").Value = -75761877
$ws.Range("H116").Value = 6934.048
$ws.Range("I116").Value = 1156.7646
$ws.Range("J116").Value = 31487.5
$ws.Range("K116").Value = 1156.7646
$ws.Range("L116").Value = 31487.5
$ws.Range("M116").Value = 1137.2354
$ws.Range("N116").Value = -36075.5
$ws.Range("H117").Value = 26666.666
$ws.Range("J117").Value = 26666.666
$ws.Range("L117").Value = 26666.666
$ws.Range("N117").Value = -35844.666
$ws.Range("H119").Value = 67169.75
$ws.Range("J119").Value = 67169.75
$ws.Range("L119").Value = 67169.75
$ws.Range("N119").Value = -76845.75
$ws.Range("H122").Value = 6146.5835
$ws.Range("I122").Value = 7499.75
$ws.Range("J122").Value = 5470
$ws.Range("K122").Value = 22499.25
$ws.Range("L122").Value = 16410
$ws.Range("M122").Value = -20049.25
$ws.Range("N122").Value = -21310
$ws.Range("H132").Value = 2257.7568
$ws.Range("I132").Value = 1668.0358
$ws.Range("K132").Value = 5004.107400000001
$ws.Range("M132").Value = -2474.107400000001
$ws.Range("H136").Value = 4124.1304
$ws.Range("I136").Value = 3354.182
$ws.Range("J136").Value = 4829.9165
$ws.Range("K136").Value = 10062.546
$ws.Range("L136").Value = 14489.7495
$ws.Range("M136").Value = -7512.545999999998
$ws.Range("N136").Value = -19589.7495

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6934.048
$ws.Range("I3").Value = 1156.7646
$ws.Range("J3").Value = 31487.5
$ws.Range("K3").Value = 1156.7646
$ws.Range("L3").Value = 31487.5
$ws.Range("M3").Value = -1042.7646
$ws.Range("N3").Value = -31715.5
$ws.Range("H99").Value = 4742.5557
$ws.Range("I99").Value = 4097.5713
$ws.Range("K99").Value = 4097.5713
$ws.Range("M99").Value = -2599.5713

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22510.846
$ws.Range("I31").Value = 1892.3658
$ws.Range("K31").Value = 1892.3658
$ws.Range("M31").Value = -1597.3658
$ws.Range("H34").Value = 22510.846
$ws.Range("I34").Value = 1892.3658
$ws.Range("K34").Value = 1892.3658
$ws.Range("M34").Value = -1690.3658
$ws.Range("H35").Value = 1964.2
$ws.Range("I35").Value = 1705.25
$ws.Range("K35").Value = 1705.25
$ws.Range("M35").Value = -1411.25
$ws.Range("H58").Value = 3635.5
$ws.Range("I58").Value = 1629.3125
$ws.Range("K58").Value = 1629.3125
$ws.Range("M58").Value = -1426.3125
$ws.Range("H105").Value = 2779.5557
$ws.Range("I105").Value = 1403.909
$ws.Range("J105").Value = 4941.2856
$ws.Range("K105").Value = 1403.909
$ws.Range("L105").Value = 4941.2856
$ws.Range("M105").Value = 343.0909999999999
$ws.Range("N105").Value = -8435.285599999999
$ws.Range("H107").Value = 1525.4
$ws.Range("I107").Value = 1245
$ws.Range("K107").Value = 1245
$ws.Range("M107").Value = 675
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -52620
$ws.Range("H122").Value = 5442.8667
$ws.Range("I122").Value = 1376.6364
$ws.Range("J122").Value = 16625
$ws.Range("K122").Value = 4129.9092
$ws.Range("L122").Value = 49875
$ws.Range("M122").Value = -1679.9092
$ws.Range("N122").Value = -54775
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H136").Value = 3635.5
$ws.Range("I136").Value = 1629.3125
$ws.Range("K136").Value = 4887.9375
$ws.Range("M136").Value = -2337.9375
$ws.Range("H138").Value = 40897.6
$ws.Range("J138").Value = 40897.6
$ws.Range("L138").Value = 40897.6
$ws.Range("N138").Value = -51177.6
$ws.Range("H141").Value = 223379.36
$ws.Range("J141").Value = 261908.22
$ws.Range("L141").Value = 261908.22
$ws.Range("N141").Value = -272268.22

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6798445
$ws.Range("J131").Value = 4420998
$ws.Range("L131").Value = 13262994
$ws.Range("N131").Value = -13273074
$ws.Range("H140").Value = 2465.4348
$ws.Range("I140").Value = 1641.1578
$ws.Range("J140").Value = 6380.75
$ws.Range("K140").Value = 4923.4734
$ws.Range("L140").Value = 19142.25
$ws.Range("M140").Value = 256.5266000000001
$ws.Range("N140").Value = -29502.25

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 99999.5
$ws.Range("J124").Value = 99999.5
$ws.Range("L124").Value = 99999.5
$ws.Range("N124").Value = -109819.5
$ws.Range("H132").Value = 2096.0981
$ws.Range("I132").Value = 1774.1957
$ws.Range("K132").Value = 5322.5871
$ws.Range("M132").Value = -2792.5871

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 28596092
$ws.Range("I64").Value = 40019284
$ws.Range("J64").Value = 38114
$ws.Range("K64").Value = 40019284
$ws.Range("L64").Value = 38114
$ws.Range("M64").Value = -40019036
$ws.Range("N64").Value = -38610
$ws.Range("H67").Value = 28596092
$ws.Range("I67").Value = 40019284
$ws.Range("J67").Value = 38114
$ws.Range("K67").Value = 40019284
$ws.Range("L67").Value = 38114
$ws.Range("M67").Value = -40018426
$ws.Range("N67").Value = -39830
$ws.Range("H107").Value = 1095.5
$ws.Range("I107").Value = 984.7
$ws.Range("K107").Value = 2954.1
$ws.Range("M107").Value = -1034.1
$ws.Range("H136").Value = 3119.2
$ws.Range("I136").Value = 1809.8485
$ws.Range("K136").Value = 5429.5455
$ws.Range("M136").Value = -2879.5455
